$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values must stay as literal text even when they look numeric,
# so force text number-format before assigning, then restore default style so
# the cell keeps the plain (unstyled) look it had originally.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.520.23"
$ws.Range("E2").Value = "  -3.14%  "
Set-TextValue $ws.Range("D3") "3.269.17"
$ws.Range("E3").Value = "  -5.47%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue $ws.Range("D5") "592.75"
$ws.Range("E5").Value = "  -2.96%  "
Set-TextValue $ws.Range("D6") "150.93"
$ws.Range("E6").Value = "  -9.77%  "
$ws.Range("E7").Value = "  -0.01%  "
Set-TextValue $ws.Range("D8") "3.261.09"
$ws.Range("E8").Value = "  -5.74%  "
$ws.Range("E9").Value = "  -8.45%  "
$ws.Range("E10").Value = "  -10.41%  "
$ws.Range("E11").Value = "  -4.68%  "
Set-TextValue $ws.Range("D12") "0.508"
$ws.Range("E12").Value = "  -9.99%  "
Set-TextValue $ws.Range("D13") "38.65"
$ws.Range("E13").Value = "  -12.98%  "
Set-TextValue $ws.Range("D14") "0.0000247"
$ws.Range("E14").Value = "  -8.27%  "
Set-TextValue $ws.Range("D15") "3.791.03"
$ws.Range("E15").Value = "  -5.59%  "
Set-TextValue $ws.Range("D16") "67.501.62"
$ws.Range("E16").Value = "  -3.20%  "
Set-TextValue $ws.Range("D17") "3.269.05"
$ws.Range("E17").Value = "  -5.39%  "
$ws.Range("E18").Value = "  -5.18%  "
Set-TextValue $ws.Range("D19") "534.14"
$ws.Range("E19").Value = "  -8.71%  "
$ws.Range("E20").Value = "  -12.41%  "
Set-TextValue $ws.Range("D21") "15.06"
$ws.Range("E21").Value = "  -12.24%  "
$ws.Range("E22").Value = "  -10.74%  "
$ws.Range("E23").Value = "  -11.67%  "
Set-TextValue $ws.Range("D24") "85.77"
$ws.Range("E24").Value = "  -9.97%  "
$ws.Range("E25").Value = "  -10.23%  "
Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -10.34%  "
Set-TextValue $ws.Range("D28") "8.14"
$ws.Range("E28").Value = "  -6.02%  "
$ws.Range("E29").Value = "  -12.08%  "
Set-TextValue $ws.Range("D30") "29.31"
$ws.Range("E30").Value = "  -10.95%  "
Set-TextValue $ws.Range("D31") "2.71"
$ws.Range("E31").Value = "  -4.08%  "
Set-TextValue $ws.Range("D32") "1.17"
$ws.Range("E32").Value = "  -6.18%  "
$ws.Range("E33").Value = "  -15.27%  "
Set-TextValue $ws.Range("D34") "5.77"
$ws.Range("E34").Value = "  -12.30%  "
Set-TextValue $ws.Range("D35") "524.20"
$ws.Range("E35").Value = "  -9.23%  "
$ws.Range("E36").Value = "  -0.05%  "
Set-TextValue $ws.Range("D37") "0.0458"
$ws.Range("E37").Value = "  -5.00%  "
Set-TextValue $ws.Range("D38") "53.58"
$ws.Range("E38").Value = "  -4.53%  "
$ws.Range("E39").Value = "  -10.43%  "
$ws.Range("E40").Value = "  -14.89%  "
Set-TextValue $ws.Range("D42") "2.83"
$ws.Range("E42").Value = "  -9.69%  "
Set-TextValue $ws.Range("D43") "2.946.01"
$ws.Range("E43").Value = "  -9.34%  "
Set-TextValue $ws.Range("D44") "0.268"
$ws.Range("E44").Value = "  -9.44%  "
$ws.Range("E45").Value = "  -14.89%  "
Set-TextValue $ws.Range("D46") "2.19"
$ws.Range("E46").Value = "  -8.65%  "
Set-TextValue $ws.Range("D47") "26.87"
$ws.Range("E47").Value = "  -12.58%  "
$ws.Range("E48").Value = "  -0.11%  "
Set-TextValue $ws.Range("D49") "2.34"
$ws.Range("E49").Value = "  -15.18%  "
$ws.Range("E50").Value = "  -9.38%  "
Set-TextValue $ws.Range("D51") "123.50"
$ws.Range("E51").Value = "  -7.30%  "
